$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 174 updates
$ws.Range("B174").Value = 12071
$ws.Range("C174").Value = 19129
$ws.Range("D174").Value = 67225
$ws.Range("E174").Value = 48096
$ws.Range("F174").Value = 62522
$ws.Range("J174").Value = -5888
$ws.Range("K174").Value = 10870
$ws.Range("L174").Value = 16758
$ws.Range("M174").Value = 8186
$ws.Range("S174").Value = 1015
$ws.Range("U174").Value = 1212
$ws.Range("V174").Value = -796
$ws.Range("X174").Value = 3396
$ws.Range("Y174").Value = -7057
$ws.Range("Z174").Value = 57283

# Row 175 updates
$ws.Range("B175").Value = 13256
$ws.Range("C175").Value = 20016
$ws.Range("D175").Value = 63462
$ws.Range("E175").Value = 43447
$ws.Range("F175").Value = 59078
$ws.Range("J175").Value = -5407
$ws.Range("K175").Value = 11096
$ws.Range("M175").Value = 8164
$ws.Range("S175").Value = 796
$ws.Range("U175").Value = 1287
$ws.Range("V175").Value = -468
$ws.Range("X175").Value = 3322
$ws.Range("Y175").Value = -6760
$ws.Range("Z175").Value = 56655

# Row 176 updates
$ws.Range("B176").Value = 13129
$ws.Range("C176").Value = 20354
$ws.Range("D176").Value = 60903
$ws.Range("E176").Value = 40549
$ws.Range("F176").Value = 56397
$ws.Range("J176").Value = -4545
$ws.Range("K176").Value = 11202
$ws.Range("L176").Value = 15746
$ws.Range("S176").Value = 800
$ws.Range("T176").Value = 2122
$ws.Range("U176").Value = 1322
$ws.Range("V176").Value = -133
$ws.Range("X176").Value = 2696
$ws.Range("Y176").Value = -7225
$ws.Range("Z176").Value = 56858

# Row 177 updates (B-F existing values changed) and new values G-X
$ws.Range("B177").Value = 13635
$ws.Range("C177").Value = 19817
$ws.Range("D177").Value = 59465
$ws.Range("E177").Value = 39648
$ws.Range("F177").Value = 56850
$ws.Range("G177").Value = 14674
$ws.Range("H177").Value = 32945
$ws.Range("I177").Value = 18271
$ws.Range("J177").Value = -3989
$ws.Range("K177").Value = 11809
$ws.Range("L177").Value = 15798
$ws.Range("M177").Value = 8308
$ws.Range("N177").Value = 8569
$ws.Range("O177").Value = 261
$ws.Range("P177").Value = 40
$ws.Range("Q177").Value = 1332
$ws.Range("R177").Value = 1291
$ws.Range("S177").Value = 974
$ws.Range("T177").Value = 2229
$ws.Range("U177").Value = 1255
$ws.Range("V177").Value = -189
$ws.Range("W177").Value = 2581
$ws.Range("X177").Value = 2770
$ws.Range("Y177").Value = -6183
$ws.Range("Z177").Value = 60503

$wb.Save()
